# Insert a new row 4 ("climate_change_factor_gnrl_hydropower_availability")
# into the "strategy_id-0" sheet of the Colombia calibrated model-input
# workbook, shifting the previous rows 4-11 (elasticity_gnrl_rate_occupancy_to_gdppc
# ... population_gnrl_urban) down to rows 5-12, per the commit message
# ("updated sampling ranges for experiment and added climate change factor
# to hydropower").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

# Shift existing data rows (old row 4 onward) down by one row.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new variable.
$ws.Range("A4").Value = "General"
$ws.Range("B4").Value = "climate_change_factor_gnrl_hydropower_availability"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.5

$cols = @("J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z", `
          "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS")
foreach ($col in $cols) {
    $ws.Range($col + "4").Value = 1
}
